$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.955.57'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '1.642.79'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.41'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5087'
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2568'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06397'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07768'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.306'
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').Value = '1.642.64'
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5456'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').Value = '0.0₅7849'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.80'
$ws.Range('D17').Value = '25.989.55'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.005'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '197.95'
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.440'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.970'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.043'
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.877'
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.92'
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1146'
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.893'
$ws.Range('E27').Value = '  +2.45%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.237'
$ws.Range('E29').Value = '  -0.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05031'
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.265'
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.191'
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.544'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.8956'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('E36').Value = '  -1.89%  '
$ws.Range('D37').Value = '1.128.35'
$ws.Range('E37').Value = '  -3.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5501'
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01554'
$ws.Range('E39').Value = '  -0.57%  '
$ws.Range('B40').Value = 'BabyDogeCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D40').Value = '0.0₈131'
$ws.Range('E40').Value = '  +19.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.553'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.005'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.636'
$ws.Range('E43').Value = '  -1.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8175'
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.94'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '1.777.54'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4533'
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.92'
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05081'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.003'
$ws.Range('E51').Value = '  -0.24%  '
